# Trading update: 2026-02-18 10:42:07
#
# A new MarketMaking trade (#37, closing-state snapshot at 10:42:03) arrived.
# The previously-tracked OPEN trades (#32-#36) on "All Trades" lose their
# "live position" columns (Capital After / slippages / confidence / entry
# reason / duration) since they are no longer the most-recent open trade for
# their strategy, and the brand new trade (#37) is appended with those
# columns populated. The "MarketMaking" strategy snapshot sheet is trimmed
# down to just the newest open trade.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "All Trades": clear the live-position columns (G, K:Q) on the rows that
#    are no longer the latest open MarketMaking trade, then append the new
#    trade as row 38.
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

foreach ($r in 33..37) {
    $allTrades.Cells.Item($r, 7).Value = 0        # G: Exit Price
    $allTrades.Cells.Item($r, 11).ClearContents()  # K: Capital After
    $allTrades.Cells.Item($r, 12).ClearContents()  # L: Entry Slippage (bps)
    $allTrades.Cells.Item($r, 13).ClearContents()  # M: Exit Slippage (bps)
    $allTrades.Cells.Item($r, 14).ClearContents()  # N: Confidence
    $allTrades.Cells.Item($r, 15).ClearContents()  # O: Entry Reason
    $allTrades.Cells.Item($r, 17).ClearContents()  # Q: Duration (min)
}

$newRow = 38
$allTrades.Cells.Item($newRow, 1).Value = 37
$allTrades.Cells.Item($newRow, 2).Value = "'2026-02-18"
$allTrades.Cells.Item($newRow, 3).Value = "'10:42:03"
$allTrades.Cells.Item($newRow, 4).Value = "MarketMaking"
$allTrades.Cells.Item($newRow, 5).Value = "UP"
$allTrades.Cells.Item($newRow, 6).Value = 0.42
$allTrades.Cells.Item($newRow, 8).Value = "OPEN"
$allTrades.Cells.Item($newRow, 9).Value = 0
$allTrades.Cells.Item($newRow, 10).Value = 0
$allTrades.Cells.Item($newRow, 11).Value = 100
$allTrades.Cells.Item($newRow, 12).Value = 0
$allTrades.Cells.Item($newRow, 13).Value = 0
$allTrades.Cells.Item($newRow, 14).Value = 0.6
$allTrades.Cells.Item($newRow, 15).Value = "Normal spread capture: 198 bps"
$allTrades.Cells.Item($newRow, 17).Value = 0

# ---------------------------------------------------------------------------
# 2) "MarketMaking" strategy sheet: drop the now-stale open trades (rows
#    3-6) and refresh the remaining row with the new trade's details.
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("A3:Q6").EntireRow.Delete()

$mm.Cells.Item(2, 1).Value = 37
$mm.Cells.Item(2, 3).Value = "'10:42:03"
$mm.Cells.Item(2, 6).Value = 0.42
$mm.Cells.Item(2, 15).Value = "Normal spread capture: 198 bps"
